$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # LP1912
$ws2 = $wb.Worksheets.Item(2)   # LP1912-215
$ws3 = $wb.Worksheets.Item(3)   # 6203-6173

# ---- Update header summary cells (A2: last update time, A3: total rows) ----
$ws1.Range("A2").Value = "Última actualización: 13:54:15"
$ws1.Range("A3").Value = "Total filas: 281"

$ws2.Range("A2").Value = "Última actualización: 13:54:15"
$ws2.Range("A3").Value = "Total filas: 74"

$ws3.Range("A2").Value = "Última actualización: 13:54:15"
$ws3.Range("A3").Value = "Total filas: 281"

# ---- Sheet 1 (LP1912) data changes ----
$ws1.Range("C67").Value = "14_ABASTO"
$ws1.Range("C68").Value = "215A_LA PLATA"
$ws1.Range("C146").Value = "16_SANTA ANA"
$ws1.Range("C147").Value = "11_ETCHEVERRY"
$ws1.Range("A182").Value = "10:35:49"
$ws1.Range("C182").Value = "215B_EL PATO"
$ws1.Range("D182").Value = 69
$ws1.Range("A183").Value = "09:57:03"
$ws1.Range("C183").Value = "215C_LA PLATA"
$ws1.Range("D183").Value = 107
$ws1.Range("A191").Value = "11:58:34"
$ws1.Range("C191").Value = "23_HERNANDEZ"
$ws1.Range("D191").Value = 7
$ws1.Range("A192").Value = "10:35:49"
$ws1.Range("C192").Value = "17_ROMERO"
$ws1.Range("D192").Value = 90
$ws1.Range("A200").Value = "12:29:23"
$ws1.Range("C200").Value = "11_ETCHEVERRY"
$ws1.Range("D200").Value = 0
$ws1.Range("A201").Value = "10:51:31"
$ws1.Range("C201").Value = "215C_EL PATO"
$ws1.Range("D201").Value = 98
$ws1.Range("A216").Value = "12:58:39"
$ws1.Range("C216").Value = "10_OLMOS"
$ws1.Range("D216").Value = 3
$ws1.Range("E216").Value = "LP1912"
$ws1.Range("A218").Value = "11:45:06"
$ws1.Range("C218").Value = "16_SANTA ANA"
$ws1.Range("D218").Value = 76
$ws1.Range("E218").Value = ""
$ws1.Range("A234").Value = "11:25:38"
$ws1.Range("C234").Value = "17_ROMERO"
$ws1.Range("D234").Value = 116
$ws1.Range("E234").Value = ""
$ws1.Range("A235").Value = "12:58:39"
$ws1.Range("C235").Value = "16_SANTA ANA"
$ws1.Range("D235").Value = 23
$ws1.Range("E235").Value = "LP1912"
$ws1.Range("A251").Value = "13:54:15"
$ws1.Range("B251").Value = "13:55"
$ws1.Range("C251").Value = "17_ROMERO"
$ws1.Range("D251").Value = 1
$ws1.Range("B253").Value = "14:01"
$ws1.Range("C253").Value = "16_SANTA ANA"
$ws1.Range("D253").Value = 26
$ws1.Range("A254").Value = "13:35:25"
$ws1.Range("B254").Value = "14:03"
$ws1.Range("C254").Value = "23_HERNANDEZ"
$ws1.Range("D254").Value = 28
$ws1.Range("A255").Value = "13:54:15"
$ws1.Range("B255").Value = "14:05"
$ws1.Range("C255").Value = "14_ABASTO"
$ws1.Range("D255").Value = 11
$ws1.Range("B256").Value = "14:11"
$ws1.Range("C256").Value = "15_ABASTO"
$ws1.Range("D256").Value = 73
$ws1.Range("A257").Value = "13:54:15"
$ws1.Range("B257").Value = "14:13"
$ws1.Range("C257").Value = "16_SANTA ANA"
$ws1.Range("D257").Value = 19
$ws1.Range("E257").Value = "LP1912"
$ws1.Range("A258").Value = "13:35:25"
$ws1.Range("B258").Value = "14:14"
$ws1.Range("C258").Value = "10_OLMOS"
$ws1.Range("D258").Value = 39
$ws1.Range("A259").Value = "12:58:39"
$ws1.Range("B259").Value = "14:16"
$ws1.Range("C259").Value = "27_EL RETIRO"
$ws1.Range("D259").Value = 78
$ws1.Range("B260").Value = "14:17"
$ws1.Range("C260").Value = "27_EL RETIRO"
$ws1.Range("D260").Value = 108
$ws1.Range("A261").Value = "12:58:39"
$ws1.Range("B261").Value = "14:20"
$ws1.Range("C261").Value = "215C_EL PATO"
$ws1.Range("D261").Value = 82
$ws1.Range("A262").Value = "12:38:18"
$ws1.Range("B262").Value = "14:21"
$ws1.Range("C262").Value = "215C_EL PATO"
$ws1.Range("D262").Value = 103
$ws1.Range("A263").Value = "12:29:23"
$ws1.Range("B263").Value = "14:24"
$ws1.Range("C263").Value = "11_ETCHEVERRY"
$ws1.Range("D263").Value = 115
$ws1.Range("E263").Value = ""
$ws1.Range("A264").Value = "12:41:18"
$ws1.Range("B264").Value = "14:25"
$ws1.Range("C264").Value = "11_ETCHEVERRY"
$ws1.Range("D264").Value = 104
$ws1.Range("A265").Value = "13:54:15"
$ws1.Range("B265").Value = "14:31"
$ws1.Range("C265").Value = "15_ABASTO"
$ws1.Range("D265").Value = 37
$ws1.Range("A266").Value = "12:58:39"
$ws1.Range("B266").Value = "14:33"
$ws1.Range("C266").Value = "215C_LA PLATA"
$ws1.Range("D266").Value = 95
$ws1.Range("A267").Value = "13:35:25"
$ws1.Range("B267").Value = "14:34"
$ws1.Range("C267").Value = "10_OLMOS"
$ws1.Range("D267").Value = 59
$ws1.Range("A268").Value = "12:46:01"
$ws1.Range("B268").Value = "14:34"
$ws1.Range("C268").Value = "215C_LA PLATA"
$ws1.Range("D268").Value = 108
$ws1.Range("A269").Value = "12:41:18"
$ws1.Range("B269").Value = "14:37"
$ws1.Range("C269").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D269").Value = 116
$ws1.Range("A270").Value = "13:35:25"
$ws1.Range("B270").Value = "14:38"
$ws1.Range("C270").Value = "23_HERNANDEZ"
$ws1.Range("D270").Value = 63
$ws1.Range("A271").Value = "12:41:18"
$ws1.Range("B271").Value = "14:40"
$ws1.Range("C271").Value = "17X38_ROMERO"
$ws1.Range("D271").Value = 119
$ws1.Range("A272").Value = "12:58:39"
$ws1.Range("B272").Value = "14:44"
$ws1.Range("C272").Value = "215B_EL PATO"
$ws1.Range("D272").Value = 106
$ws1.Range("A273").Value = "12:46:01"
$ws1.Range("B273").Value = "14:45"
$ws1.Range("C273").Value = "215B_EL PATO"
$ws1.Range("D273").Value = 119
$ws1.Range("A274").Value = "12:58:39"
$ws1.Range("B274").Value = "14:53"
$ws1.Range("D274").Value = 115
$ws1.Range("A275").Value = "12:58:39"
$ws1.Range("B275").Value = "14:53"
$ws1.Range("C275").Value = "215A_LA PLATA"
$ws1.Range("D275").Value = 115
$ws1.Range("B276").Value = "14:56"
$ws1.Range("C276").Value = "215A_LA PLATA"
$ws1.Range("D276").Value = 81
$ws1.Range("B277").Value = "15:01"
$ws1.Range("C277").Value = "81_EL PELIGRO"
$ws1.Range("D277").Value = 86
$ws1.Range("A278").Value = "13:54:15"
$ws1.Range("B278").Value = "15:02"
$ws1.Range("C278").Value = "215A_LA PLATA"
$ws1.Range("D278").Value = 68
$ws1.Range("E278").Value = "LP1912"
$ws1.Range("A279").Value = "13:35:25"
$ws1.Range("B279").Value = "15:04"
$ws1.Range("C279").Value = "14_ABASTO"
$ws1.Range("D279").Value = 89
$ws1.Range("E279").Value = "LP1912"
$ws1.Range("A280").Value = "13:54:15"
$ws1.Range("B280").Value = "15:05"
$ws1.Range("C280").Value = "14_ABASTO"
$ws1.Range("D280").Value = 71
$ws1.Range("E280").Value = "LP1912"
$ws1.Range("A281").Value = "13:35:25"
$ws1.Range("B281").Value = "15:17"
$ws1.Range("C281").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D281").Value = 102
$ws1.Range("E281").Value = "LP1912"
$ws1.Range("A282").Value = "13:35:25"
$ws1.Range("B282").Value = "15:24"
$ws1.Range("C282").Value = "11_ETCHEVERRY"
$ws1.Range("D282").Value = 109
$ws1.Range("E282").Value = "LP1912"
$ws1.Range("A283").Value = "13:35:25"
$ws1.Range("B283").Value = "15:25"
$ws1.Range("C283").Value = "215C_EL PATO"
$ws1.Range("D283").Value = 110
$ws1.Range("E283").Value = "LP1912"
$ws1.Range("A284").Value = "13:54:15"
$ws1.Range("B284").Value = "15:25"
$ws1.Range("C284").Value = "11_ETCHEVERRY"
$ws1.Range("D284").Value = 91
$ws1.Range("E284").Value = "LP1912"
$ws1.Range("A285").Value = "13:54:15"
$ws1.Range("B285").Value = "15:25"
$ws1.Range("C285").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D285").Value = 91
$ws1.Range("E285").Value = "LP1912"
$ws1.Range("A286").Value = "13:54:15"
$ws1.Range("B286").Value = "15:36"
$ws1.Range("C286").Value = "17X38_ROMERO"
$ws1.Range("D286").Value = 102
$ws1.Range("E286").Value = "LP1912"

# ---- Sheet 2 (LP1912-215) data changes ----
$ws2.Range("A50").Value = "09:57:03"
$ws2.Range("C50").Value = "215C_LA PLATA"
$ws2.Range("D50").Value = 107
$ws2.Range("A51").Value = "10:35:49"
$ws2.Range("C51").Value = "215B_EL PATO"
$ws2.Range("D51").Value = 69
$ws2.Range("A78").Value = "13:54:15"
$ws2.Range("B78").Value = "15:02"
$ws2.Range("C78").Value = "215A_LA PLATA"
$ws2.Range("D78").Value = 68
$ws2.Range("A79").Value = "13:35:25"
$ws2.Range("B79").Value = "15:25"
$ws2.Range("C79").Value = "215C_EL PATO"
$ws2.Range("D79").Value = 110
$ws2.Range("E79").Value = "LP1912"

# ---- Sheet 3 (6203-6173) data changes (mirrors sheet 1) ----
$ws3.Range("C67").Value = "14_ABASTO"
$ws3.Range("C68").Value = "215A_LA PLATA"
$ws3.Range("C146").Value = "16_SANTA ANA"
$ws3.Range("C147").Value = "11_ETCHEVERRY"
$ws3.Range("A182").Value = "10:35:49"
$ws3.Range("C182").Value = "215B_EL PATO"
$ws3.Range("D182").Value = 69
$ws3.Range("A183").Value = "09:57:03"
$ws3.Range("C183").Value = "215C_LA PLATA"
$ws3.Range("D183").Value = 107
$ws3.Range("A191").Value = "11:58:34"
$ws3.Range("C191").Value = "23_HERNANDEZ"
$ws3.Range("D191").Value = 7
$ws3.Range("A192").Value = "10:35:49"
$ws3.Range("C192").Value = "17_ROMERO"
$ws3.Range("D192").Value = 90
$ws3.Range("A200").Value = "12:29:23"
$ws3.Range("C200").Value = "11_ETCHEVERRY"
$ws3.Range("D200").Value = 0
$ws3.Range("A201").Value = "10:51:31"
$ws3.Range("C201").Value = "215C_EL PATO"
$ws3.Range("D201").Value = 98
$ws3.Range("A216").Value = "12:58:39"
$ws3.Range("C216").Value = "10_OLMOS"
$ws3.Range("D216").Value = 3
$ws3.Range("E216").Value = "LP1912"
$ws3.Range("A218").Value = "11:45:06"
$ws3.Range("C218").Value = "16_SANTA ANA"
$ws3.Range("D218").Value = 76
$ws3.Range("E218").Value = ""
$ws3.Range("A234").Value = "11:25:38"
$ws3.Range("C234").Value = "17_ROMERO"
$ws3.Range("D234").Value = 116
$ws3.Range("E234").Value = ""
$ws3.Range("A235").Value = "12:58:39"
$ws3.Range("C235").Value = "16_SANTA ANA"
$ws3.Range("D235").Value = 23
$ws3.Range("E235").Value = "LP1912"
$ws3.Range("A251").Value = "13:54:15"
$ws3.Range("B251").Value = "13:55"
$ws3.Range("C251").Value = "17_ROMERO"
$ws3.Range("D251").Value = 1
$ws3.Range("B253").Value = "14:01"
$ws3.Range("C253").Value = "16_SANTA ANA"
$ws3.Range("D253").Value = 26
$ws3.Range("A254").Value = "13:35:25"
$ws3.Range("B254").Value = "14:03"
$ws3.Range("C254").Value = "23_HERNANDEZ"
$ws3.Range("D254").Value = 28
$ws3.Range("A255").Value = "13:54:15"
$ws3.Range("B255").Value = "14:05"
$ws3.Range("C255").Value = "14_ABASTO"
$ws3.Range("D255").Value = 11
$ws3.Range("B256").Value = "14:11"
$ws3.Range("C256").Value = "15_ABASTO"
$ws3.Range("D256").Value = 73
$ws3.Range("A257").Value = "13:54:15"
$ws3.Range("B257").Value = "14:13"
$ws3.Range("C257").Value = "16_SANTA ANA"
$ws3.Range("D257").Value = 19
$ws3.Range("E257").Value = "LP1912"
$ws3.Range("A258").Value = "13:35:25"
$ws3.Range("B258").Value = "14:14"
$ws3.Range("C258").Value = "10_OLMOS"
$ws3.Range("D258").Value = 39
$ws3.Range("A259").Value = "12:58:39"
$ws3.Range("B259").Value = "14:16"
$ws3.Range("C259").Value = "27_EL RETIRO"
$ws3.Range("D259").Value = 78
$ws3.Range("B260").Value = "14:17"
$ws3.Range("C260").Value = "27_EL RETIRO"
$ws3.Range("D260").Value = 108
$ws3.Range("A261").Value = "12:58:39"
$ws3.Range("B261").Value = "14:20"
$ws3.Range("C261").Value = "215C_EL PATO"
$ws3.Range("D261").Value = 82
$ws3.Range("A262").Value = "12:38:18"
$ws3.Range("B262").Value = "14:21"
$ws3.Range("C262").Value = "215C_EL PATO"
$ws3.Range("D262").Value = 103
$ws3.Range("A263").Value = "12:29:23"
$ws3.Range("B263").Value = "14:24"
$ws3.Range("C263").Value = "11_ETCHEVERRY"
$ws3.Range("D263").Value = 115
$ws3.Range("E263").Value = ""
$ws3.Range("A264").Value = "12:41:18"
$ws3.Range("B264").Value = "14:25"
$ws3.Range("C264").Value = "11_ETCHEVERRY"
$ws3.Range("D264").Value = 104
$ws3.Range("A265").Value = "13:54:15"
$ws3.Range("B265").Value = "14:31"
$ws3.Range("C265").Value = "15_ABASTO"
$ws3.Range("D265").Value = 37
$ws3.Range("A266").Value = "12:58:39"
$ws3.Range("B266").Value = "14:33"
$ws3.Range("C266").Value = "215C_LA PLATA"
$ws3.Range("D266").Value = 95
$ws3.Range("A267").Value = "13:35:25"
$ws3.Range("B267").Value = "14:34"
$ws3.Range("C267").Value = "10_OLMOS"
$ws3.Range("D267").Value = 59
$ws3.Range("A268").Value = "12:46:01"
$ws3.Range("B268").Value = "14:34"
$ws3.Range("C268").Value = "215C_LA PLATA"
$ws3.Range("D268").Value = 108
$ws3.Range("A269").Value = "12:41:18"
$ws3.Range("B269").Value = "14:37"
$ws3.Range("C269").Value = "16_P MOR-SANTA ANA"
$ws3.Range("D269").Value = 116
$ws3.Range("A270").Value = "13:35:25"
$ws3.Range("B270").Value = "14:38"
$ws3.Range("C270").Value = "23_HERNANDEZ"
$ws3.Range("D270").Value = 63
$ws3.Range("A271").Value = "12:41:18"
$ws3.Range("B271").Value = "14:40"
$ws3.Range("C271").Value = "17X38_ROMERO"
$ws3.Range("D271").Value = 119
$ws3.Range("A272").Value = "12:58:39"
$ws3.Range("B272").Value = "14:44"
$ws3.Range("C272").Value = "215B_EL PATO"
$ws3.Range("D272").Value = 106
$ws3.Range("A273").Value = "12:46:01"
$ws3.Range("B273").Value = "14:45"
$ws3.Range("C273").Value = "215B_EL PATO"
$ws3.Range("D273").Value = 119
$ws3.Range("A274").Value = "12:58:39"
$ws3.Range("B274").Value = "14:53"
$ws3.Range("D274").Value = 115
$ws3.Range("A275").Value = "12:58:39"
$ws3.Range("B275").Value = "14:53"
$ws3.Range("C275").Value = "215A_LA PLATA"
$ws3.Range("D275").Value = 115
$ws3.Range("B276").Value = "14:56"
$ws3.Range("C276").Value = "215A_LA PLATA"
$ws3.Range("D276").Value = 81
$ws3.Range("B277").Value = "15:01"
$ws3.Range("C277").Value = "81_EL PELIGRO"
$ws3.Range("D277").Value = 86
$ws3.Range("A278").Value = "13:54:15"
$ws3.Range("B278").Value = "15:02"
$ws3.Range("C278").Value = "215A_LA PLATA"
$ws3.Range("D278").Value = 68
$ws3.Range("E278").Value = "LP1912"
$ws3.Range("A279").Value = "13:35:25"
$ws3.Range("B279").Value = "15:04"
$ws3.Range("C279").Value = "14_ABASTO"
$ws3.Range("D279").Value = 89
$ws3.Range("E279").Value = "LP1912"
$ws3.Range("A280").Value = "13:54:15"
$ws3.Range("B280").Value = "15:05"
$ws3.Range("C280").Value = "14_ABASTO"
$ws3.Range("D280").Value = 71
$ws3.Range("E280").Value = "LP1912"
$ws3.Range("A281").Value = "13:35:25"
$ws3.Range("B281").Value = "15:17"
$ws3.Range("C281").Value = "16_P MOR-SANTA ANA"
$ws3.Range("D281").Value = 102
$ws3.Range("E281").Value = "LP1912"
$ws3.Range("A282").Value = "13:35:25"
$ws3.Range("B282").Value = "15:24"
$ws3.Range("C282").Value = "11_ETCHEVERRY"
$ws3.Range("D282").Value = 109
$ws3.Range("E282").Value = "LP1912"
$ws3.Range("A283").Value = "13:35:25"
$ws3.Range("B283").Value = "15:25"
$ws3.Range("C283").Value = "215C_EL PATO"
$ws3.Range("D283").Value = 110
$ws3.Range("E283").Value = "LP1912"
$ws3.Range("A284").Value = "13:54:15"
$ws3.Range("B284").Value = "15:25"
$ws3.Range("C284").Value = "11_ETCHEVERRY"
$ws3.Range("D284").Value = 91
$ws3.Range("E284").Value = "LP1912"
$ws3.Range("A285").Value = "13:54:15"
$ws3.Range("B285").Value = "15:25"
$ws3.Range("C285").Value = "16_P MOR-SANTA ANA"
$ws3.Range("D285").Value = 91
$ws3.Range("E285").Value = "LP1912"
$ws3.Range("A286").Value = "13:54:15"
$ws3.Range("B286").Value = "15:36"
$ws3.Range("C286").Value = "17X38_ROMERO"
$ws3.Range("D286").Value = 102
$ws3.Range("E286").Value = "LP1912"
